# Updated symbol list on Mon Jan  9 04:55:24 UTC 2023 with GitHub Actions
# Refreshes coin prices / 1h volume percentages, and shifts the exchange-token
# rows (7-18) down by one position to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the numeric-looking value as text
# (matching the workbook's original inlineStr/text cell format).

$ws.Range("D2").Value = "'278.57"
$ws.Range("E2").Value = "'6.64%"

$ws.Range("D3").Value = "'27.27"
$ws.Range("E3").Value = "'0.84%"

$ws.Range("D4").Value = "'4.824"
$ws.Range("E4").Value = "'2.85%"

$ws.Range("D5").Value = "'0.06282"
$ws.Range("E5").Value = "'0.97%"

$ws.Range("D6").Value = "'6.856"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.275"
$ws.Range("E7").Value = "'3.12%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8775"
$ws.Range("E8").Value = "'3.00%"

$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.9475"
$ws.Range("E9").Value = "'3.80%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1455"
$ws.Range("E10").Value = "'4.09%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.05161"
$ws.Range("E11").Value = "'9.13%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07287"
$ws.Range("E12").Value = "'2.77%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03167"
$ws.Range("E13").Value = "'1.19%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09054"
$ws.Range("E14").Value = "'-0.03%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001553"
$ws.Range("E15").Value = "'1.82%"

$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006269"
$ws.Range("E16").Value = "'2.12%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005929"
$ws.Range("E17").Value = "'-3.51%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.468"
$ws.Range("E18").Value = "'0.80%"

$ws.Range("D19").Value = "'2.266"
$ws.Range("E19").Value = "'4.59%"

$ws.Range("E20").Value = "'-0.62%"

$ws.Range("D21").Value = "'0.1308"
$ws.Range("E21").Value = "'-0.18%"

$ws.Range("D22").Value = "'3.849"
$ws.Range("E22").Value = "'-5.64%"

$ws.Range("D23").Value = "'0.04323"
$ws.Range("E23").Value = "'1.88%"

$ws.Range("D24").Value = "'0.001175"
$ws.Range("E24").Value = "'-3.10%"

$ws.Range("D25").Value = "'0.004277"
$ws.Range("E25").Value = "'4.57%"

$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.14%"

$ws.Range("D27").Value = "'0.0001774"
$ws.Range("E27").Value = "'8.20%"

$ws.Range("D40").Value = "'0.04034"
$ws.Range("E40").Value = "'3.08%"

$ws.Range("D41").Value = "'0.006712"
$ws.Range("E41").Value = "'62.96%"

$ws.Range("D42").Value = "'0.1155"
$ws.Range("E42").Value = "'3.78%"

$ws.Range("D43").Value = "'0.01408"
$ws.Range("E43").Value = "'1.33%"

$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'-4.99%"

$ws.Range("D45").Value = "'0.00005183"
$ws.Range("E45").Value = "'1.20%"

$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.15%"

$ws.Range("D47").Value = "'2.310"
$ws.Range("E47").Value = "'636.54%"

$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.15%"

$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.15%"
